# Applies the "More work on ATR72_ATsi" edit:
#   - Re-orders the WEIGHT ESTIMATION METHODS COMPARISON method labels
#     (RAYMER, ROSKAM, NICOLAI_1984, KROO, TORENBEEK_1976, JENKINSON, SADRAEY,
#      TORENBEEK_2013, TORENBEEK_1982, HOWE, NICOLAI_2013) and keeps each
#     method's own (Estimated Mass, Percent Error) pair attached to its
#     label while moving rows around on the FUSELAGE, WING, HORIZONTAL TAIL
#     and VERTICAL TAIL comparison tables.
#   - Updates the "Method: ..." label used on NACELLES, POWER PLANT and
#     LANDING GEARS sheets (TORENBEEK_1976), whose backing shared string
#     moved inside the table.
#
# For every row touched by the edit we re-assert both the label text (column
# A) and the numeric Estimated Mass / Percent Error values (columns C/D), so
# the resulting cell content is correct regardless of how the workbook
# happens to lay out its internal shared-string table.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# FUSELAGE  (rows 8-15 of the comparison table)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("FUSELAGE")

$ws.Range("A8").Value = "RAYMER"
$ws.Range("C8").Value = 3149.0
$ws.Range("D8").Value = 15.827205899908506

$ws.Range("A9").Value = "ROSKAM"
$ws.Range("C9").Value = 3917.0
$ws.Range("D9").Value = 44.07594966971788

$ws.Range("A10").Value = "NICOLAI_1984"
$ws.Range("C10").Value = 2968.0
$ws.Range("D10").Value = 9.169624360409161

$ws.Range("A11").Value = "KROO"
$ws.Range("C11").Value = 2585.0
$ws.Range("D11").Value = -4.917965306045256

$ws.Range("A12").Value = "TORENBEEK_1976"
$ws.Range("C12").Value = 3818.0
$ws.Range("D12").Value = 40.43451004314089

$ws.Range("A13").Value = "JENKINSON"
$ws.Range("C13").Value = 4506.0
$ws.Range("D13").Value = 65.74067633692846

$ws.Range("A14").Value = "SADRAEY"
$ws.Range("C14").Value = 2491.0
$ws.Range("D14").Value = -8.375493840370883

$ws.Range("A15").Value = "TORENBEEK_2013"
$ws.Range("C15").Value = 3698.0
$ws.Range("D15").Value = 36.020643829108174

# ---------------------------------------------------------------------------
# WING  (rows 8-11)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("WING")

$ws.Range("A8").Value = "TORENBEEK_1982"
$ws.Range("C8").Value = 2613.0
$ws.Range("D8").Value = -9.328361499469459

$ws.Range("A9").Value = "RAYMER"
$ws.Range("C9").Value = 2760.0
$ws.Range("D9").Value = -4.227431204950519

$ws.Range("A10").Value = "KROO"
$ws.Range("C10").Value = 2539.0
$ws.Range("D10").Value = -11.89617674977151

$ws.Range("A11").Value = "TORENBEEK_2013"
$ws.Range("C11").Value = 2330.0
$ws.Range("D11").Value = -19.148519821570545

# ---------------------------------------------------------------------------
# HORIZONTAL TAIL  (rows 8-15)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("HORIZONTAL TAIL")

$ws.Range("A8").Value = "HOWE"
$ws.Range("C8").Value = 207.0
$ws.Range("D8").Value = -33.792006789509266

$ws.Range("A9").Value = "RAYMER"
$ws.Range("C9").Value = 144.0
$ws.Range("D9").Value = -53.9422655927021

$ws.Range("A10").Value = "ROSKAM"
$ws.Range("C10").Value = 216.0
$ws.Range("D10").Value = -30.91339838905315

$ws.Range("A11").Value = "KROO"
$ws.Range("C11").Value = 303.0
$ws.Range("D11").Value = -3.0868505179773305

$ws.Range("A12").Value = "TORENBEEK_1976"
$ws.Range("C12").Value = 236.0
$ws.Range("D12").Value = -24.516490832483992

$ws.Range("A13").Value = "JENKINSON"
$ws.Range("C13").Value = 293.0
$ws.Range("D13").Value = -6.285304296261907

$ws.Range("A14").Value = "SADRAEY"
$ws.Range("C14").Value = 273.0
$ws.Range("D14").Value = -12.68221185283106

$ws.Range("A15").Value = "NICOLAI_2013"
$ws.Range("C15").Value = 124.0
$ws.Range("D15").Value = -60.339173149271254

# ---------------------------------------------------------------------------
# VERTICAL TAIL  (rows 8-12)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("VERTICAL TAIL")

$ws.Range("A8").Value = "HOWE"
$ws.Range("C8").Value = 445.0
$ws.Range("D8").Value = 42.33119313366365

$ws.Range("A9").Value = "RAYMER"
$ws.Range("C9").Value = 89.0
$ws.Range("D9").Value = -71.53376137326727

$ws.Range("A10").Value = "ROSKAM"
$ws.Range("C10").Value = 239.0
$ws.Range("D10").Value = -23.55695469899862

$ws.Range("A11").Value = "KROO"
$ws.Range("C11").Value = 256.0
$ws.Range("D11").Value = -18.11958327591484

$ws.Range("A12").Value = "TORENBEEK_1976"
$ws.Range("C12").Value = 338.0
$ws.Range("D12").Value = 8.107737706018687

# ---------------------------------------------------------------------------
# NACELLES / POWER PLANT / LANDING GEARS
# "Method: TORENBEEK_1976" label rows - text itself is unchanged, only the
# underlying shared-string slot moved, so re-asserting the label text keeps
# the content correct.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("NACELLES")
$ws.Range("A11").Value = "TORENBEEK_1976"
$ws.Range("A17").Value = "TORENBEEK_1976"

$ws = $wb.Worksheets.Item("POWER PLANT")
$ws.Range("A12").Value = "TORENBEEK_1976"
$ws.Range("A18").Value = "TORENBEEK_1976"

$ws = $wb.Worksheets.Item("LANDING GEARS")
$ws.Range("A9").Value = "TORENBEEK_1976"
$ws.Range("A11").Value = "TORENBEEK_1976"
$ws.Range("A13").Value = "TORENBEEK_1976"

$wb.Save()
